$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"
$wsElem.Columns.Item(26).ColumnWidth = 50.33
